{"js": "// Word JS API (Office.js) edit script.\n//\n// The original document had several bullet paragraphs whose sentence-\n// starting verb was split across two (or more) separate <w:r> runs, e.g.\n//   \"Contribut\" + \"ed\" + \" to various ... biomaterials\"\n// with no formatting difference between the pieces. This script merges\n// each such run-group back into a single run by replacing the full\n// (run-spanning) text range with the same text in one `insertText` call,\n// which Office.js always writes out as a single run.\n//\n// It also introduces a new bold run around \"300-person\" inside the\n// \"Built and managed multiple CI/CD pipelines ...\" bullet, splitting the\n// previously single \" for a 300-person engineering organization...\" run\n// into \" for a \" + bold(\"300-person\") + \" engineering organization...\".\n\nasync function mergeRun(context, fullText) {\n  const results = context.document.body.search(fullText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"mergeRun: text not found -> \" + fullText);\n  }\n  // Re-writing the exact same text over the whole (possibly multi-run)\n  // range collapses it down to a single run sharing one <w:rPr/>.\n  results.items[0].insertText(fullText, \"Replace\");\n  await context.sync();\n}\n\n// 1) \"Contribut\" + \"ed\" + \" to various ...\" -> one run.\nawait mergeRun(\n  context,\n  \"Contributed to various company-wide quality initiatives for a cloud-based software company that develops software platforms powering breakthrough research on biotherapeutics, biofuels, and biomaterials\"\n);\n\n// 2) \"Enhanc\" + \"ed\" + \" test infrastructure by developing and maintaining \" -> one run.\nawait mergeRun(\n  context,\n  \"Enhanced test infrastructure by developing and maintaining \"\n);\n\n// 3) \"Built\" + \" and mana\" + \"ged\" + \" multiple \" -> one run.\nawait mergeRun(context, \"Built and managed multiple \");\n\n// 4) \"Support\" + \"ed\" + \" TechOps (Technical Operations) by creating and analyzing o\" -> one run.\nawait mergeRun(\n  context,\n  \"Supported TechOps (Technical Operations) by creating and analyzing o\"\n);\n\n// 5) Bold just \"300-person\" inside the \"... for a 300-person engineering\n//    organization ...\" run, splitting it into three runs.\nconst personResults = context.document.body.search(\"300-person\", {\n  matchCase: true,\n});\npersonResults.load(\"items\");\nawait context.sync();\nif (personResults.items.length === 0) {\n  throw new Error(\"300-person text not found\");\n}\npersonResults.items[0].font.bold = true;\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n#\n# The original document had several bullet paragraphs whose sentence-\n# starting verb was split across two (or more) separate runs, e.g.\n#   \"Contribut\" + \"ed\" + \" to various ... biomaterials\"\n# with no formatting difference between the pieces. This script merges\n# each such run-group back into a single run by running a Find/Replace\n# (wdReplaceAll) over the full sentence text - Word writes the replaced\n# range back out as a single run.\n#\n# It also introduces a new bold run around \"300-person\" inside the\n# \"Built and managed multiple CI/CD pipelines ...\" bullet, splitting the\n# previously single \" for a 300-person engineering organization...\" run\n# into \" for a \" + bold(\"300-person\") + \" engineering organization...\".\n\n$d = $word.ActiveDocument\n\nfunction Merge-Run([string]$text) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Execute(\n        $text,   # FindText\n        $false,  # MatchCase\n        $false,  # MatchWholeWord\n        $false,  # MatchWildcards\n        $false,  # MatchSoundsLike\n        $false,  # MatchAllWordForms\n        $true,   # Forward\n        1,       # Wrap (wdFindContinue)\n        $false,  # Format\n        $text,   # ReplaceWith\n        2        # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n\n# 1) \"Contribut\" + \"ed\" + \" to various ...\" -> one run.\nMerge-Run \"Contributed to various company-wide quality initiatives for a cloud-based software company that develops software platforms powering breakthrough research on biotherapeutics, biofuels, and biomaterials\"\n\n# 2) \"Enhanc\" + \"ed\" + \" test infrastructure by developing and maintaining \" -> one run.\nMerge-Run \"Enhanced test infrastructure by developing and maintaining \"\n\n# 3) \"Built\" + \" and mana\" + \"ged\" + \" multiple \" -> one run.\nMerge-Run \"Built and managed multiple \"\n\n# 4) \"Support\" + \"ed\" + \" TechOps (Technical Operations) by creating and analyzing o\" -> one run.\nMerge-Run \"Supported TechOps (Technical Operations) by creating and analyzing o\"\n\n# 5) Bold just \"300-person\" inside the \"... for a 300-person engineering\n#    organization ...\" run, splitting it into three runs.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"300-person\"\n$find.Execute() | Out-Null\n$rng.Bold = 1\n"}
